# Add the new row of data to Sheet1: A2 = "updated file"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "updated file"

# Leave the new cell selected, mirroring the author's final selection state
$ws.Range("A2").Select()
